# Refresh the "想去人数" (interest count) figures in column F that
# changed on this data pull, mirroring each exhibition/local-life row's
# update into the aggregate "全部类型" sheet as well.

$wb = $excel.ActiveWorkbook

$sheetExpo = $wb.Worksheets.Item("展览")      # exhibitions
$sheetLocal = $wb.Worksheets.Item("本地生活")  # local life
$sheetAll = $wb.Worksheets.Item("全部类型")    # all types (aggregate)

# --- 展览 (sheet1) updates, mirrored on 全部类型 ---
$sheetExpo.Range("F2").Value = 1245
$sheetExpo.Range("F5").Value = 5515
$sheetExpo.Range("F6").Value = 1762
$sheetExpo.Range("F7").Value = 6321
$sheetExpo.Range("F9").Value = 1900
$sheetExpo.Range("F10").Value = 508
$sheetExpo.Range("F16").Value = 7830
$sheetExpo.Range("F17").Value = 7830
$sheetExpo.Range("F29").Value = 1711
$sheetExpo.Range("F31").Value = 361

$sheetAll.Range("F5").Value = 1245
$sheetAll.Range("F10").Value = 5515
$sheetAll.Range("F12").Value = 1762
$sheetAll.Range("F13").Value = 6321
$sheetAll.Range("F15").Value = 1900
$sheetAll.Range("F17").Value = 508
$sheetAll.Range("F23").Value = 7830
$sheetAll.Range("F24").Value = 7830
$sheetAll.Range("F35").Value = 1711
$sheetAll.Range("F38").Value = 361

# --- 本地生活 (sheet3) updates, mirrored on 全部类型 ---
$sheetLocal.Range("F2").Value = 9537
$sheetLocal.Range("F4").Value = 675

$sheetAll.Range("F2").Value = 9537
$sheetAll.Range("F4").Value = 675
